$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = -0.4344319951986923
$ws.Range("J2").Value = 0.2341496727038916
$ws.Range("K2").Value = -0.114760828551843
$ws.Range("L2").Value = 2.729621715297354
